{"js": "// Add a reviewer comment \"cris\" anchored on the heading text\n// \"Scope and Delimitations\" (inside \"IV.  Scope and Delimitations\"),\n// matching the author's intent of tagging \"cris\" on that section.\n\nconst body = context.document.body;\n\n// Find the exact run of text the comment should be anchored to.\nconst results = body.search(\"Scope and Delimitations\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Scope and Delimitations\" in the document body.');\n}\n\n// Anchor the new comment on the found range and give it the comment text.\nconst targetRange = results.items[0];\ntargetRange.insertComment(\"cris\");\n\nawait context.sync();\n", "ps1": "# Add a reviewer comment \"cris\" anchored on the heading text\n# \"Scope and Delimitations\" (inside \"IV.  Scope and Delimitations\"),\n# matching the author's intent of tagging \"cris\" on that section.\n\n$app.UserName = \"Francis Santos\"\n$app.UserInitials = \"FS\"\n\n$d = $word.ActiveDocument\n\n# Locate the exact text the comment should be anchored to.\n$rng = $d.Content\n$rng.Find.Execute(\"Scope and Delimitations\") | Out-Null\n\n# Anchor the new comment on the found range and give it the comment text.\n$d.Comments.Add($rng, \"cris\") | Out-Null\n"}
